$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B1").EntireColumn.Delete()
$ws.Range("G1").Value = "Planilla Salud"
$ws.Range("G1").Font.Bold = $true
$ws.Columns("G").AutoFit()
